$d = $word.ActiveDocument

# 1) Merge the split "Versi" + "on" runs into a single "Version" run
#    by replacing the matched text with itself.
$d.Content.Find.Execute("Version", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Version", 2)

# 2) Change the version number from 2 to 1 (keeps it inside the
#    " 2" run, now " 1", leaving the bookmark and trailing "." run
#    untouched).
$d.Content.Find.Execute("2", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "1", 2)

# 3) Move the trailing "." (currently its own run, after the
#    _GoBack bookmark) so that it becomes part of the " 1" run,
#    producing " 1." before the bookmark and leaving nothing after it.
$periodRange = $d.Range(9, 10)
$periodRange.Delete()

$insertRange = $d.Range(9, 9)
$insertRange.InsertBefore(".")
